$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1) onto the new H1 header cell,
# then set its text, matching the other "bold/centered/bordered" header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new data cell H2 = 0 (numeric), no special style (like F2/G2).
$ws.Range("H2").Value = 0
